$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new package name "ss" as the next row in column A
$ws.Range("A53").Value = "ss"

# Match the author's saved selection state (active cell moves to the new row)
[void]$ws.Range("A53").Select()
